# Correction BR 1 / 2
# Set "Avancement" (progress) column F to "En cours Excel" for rows 53-59,
# which were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F53:F59").Value = "En cours Excel"

# Update the active selection / scroll position as reflected in the saved view.
$ws.Range("G57").Select()
